# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.813.13"
$ws.Range("E2").Value = "  -4.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.904.17"
$ws.Range("E3").Value = "  -7.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "472.25"
$ws.Range("E5").Value = "  -10.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.97"
$ws.Range("E6").Value = "  -4.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.913.45"
$ws.Range("E8").Value = "  -7.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.409"
$ws.Range("E9").Value = "  -8.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").Value = "  -7.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0983"
$ws.Range("E11").Value = "  -11.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  -14.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.123"
$ws.Range("E13").Value = "  -3.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.385.44"
$ws.Range("E14").Value = "  -8.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.55"
$ws.Range("E15").Value = "  -8.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "54.487.78"
$ws.Range("E16").Value = "  -5.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.888.77"
$ws.Range("E17").Value = "  -8.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000135"
$ws.Range("E18").Value = "  -11.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.44"
$ws.Range("E19").Value = "  -6.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.57"
$ws.Range("E20").Value = "  -11.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.15"
$ws.Range("E21").Value = "  -10.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "305.95"
$ws.Range("E22").Value = "  -11.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.445"
$ws.Range("E24").Value = "  -12.91%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "59.45"
$ws.Range("E25").Value = "  -14.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.151"
$ws.Range("E27").Value = "  -9.20%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0815"
$ws.Range("E29").Value = "  -15.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.15"
$ws.Range("E30").Value = "  -10.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.30"
$ws.Range("E31").Value = "  -9.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.11"
$ws.Range("E32").Value = "  -8.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.03"
$ws.Range("E33").Value = "  -11.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.62"
$ws.Range("E34").Value = "  -13.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "143.04"
$ws.Range("E35").Value = "  -9.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.24"
$ws.Range("E36").Value = "  -14.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.46"
$ws.Range("E37").Value = "  -12.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.22"
$ws.Range("E38").Value = "  -13.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.95"
$ws.Range("E39").Value = "  -11.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0624"
$ws.Range("E40").Value = "  -10.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.920.91"
$ws.Range("E41").Value = "  -8.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.996"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.71"
$ws.Range("E43").Value = "  -11.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.968"
$ws.Range("E44").Value = "  -10.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.608"
$ws.Range("E45").Value = "  -12.39%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.34"
$ws.Range("E46").Value = "  -8.33%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.46"
$ws.Range("E47").Value = "  -12.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.038.69"
$ws.Range("E48").Value = "  -10.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.39"
$ws.Range("E49").Value = "  -13.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0216"
$ws.Range("E50").Value = "  -7.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.93"
$ws.Range("E51").Value = "  -12.86%  "
